# Update countries & provincias Spain
#
# This workbook ("Pais" sheet) lists per-country COVID stats. The update:
#   1) Refreshes the "last updated" timestamp cell.
#   2) Refreshes several countries' statistics (totals/new/active/recovered/
#      critical/deaths-today/deaths) for this data refresh.
#   3) A handful of country rows got re-sorted in the source data, which
#      (because the row positions did NOT move) shows up as adjacent rows
#      swapping/rotating both their country name and their stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (row 1) ---
$ws.Range("A1").Value = 'Datos actualizados a 25 de Mayo de 2020 a las 21:05'

# --- Estados Unidos (row 4): stats refresh only ---
$ws.Range("B4").Value = 1697692
$ws.Range("C4").Value = 11256
$ws.Range("E4").Value = 1141461
$ws.Range("G4").Value = 336
$ws.Range("H4").Value = 99636

# --- Alemania (row 11): stats refresh only ---
$ws.Range("B11").Value = 180688
$ws.Range("C11").Value = 360
$ws.Range("E11").Value = 11075
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = 8413

# --- Rows 84-85: Bosnia y Herzegovina <-> Costa de Marfil swap ---
$ws.Range("A84").Value = 'Costa de Marfil'
$ws.Range("B84").Value = 2423
$ws.Range("C84").Value = 47
$ws.Range("D84").Value = 1257
$ws.Range("E84").Value = 1136
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 30

$ws.Range("A85").Value = 'Bosnia y Herzegovina'
$ws.Range("B85").Value = 2406
$ws.Range("C85").Value = 5
$ws.Range("D85").Value = 1696
$ws.Range("E85").Value = 564
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 146

# --- Rows 104-106: Venezuela / Libano / Guinea-Bisau rotate ---
$ws.Range("A104").Value = 'Guinea-Bisau'
$ws.Range("B104").Value = 1178
$ws.Range("C104").Value = 64
$ws.Range("D104").Value = 42
$ws.Range("E104").Value = 1129
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 7

$ws.Range("A105").Value = 'Venezuela'
$ws.Range("B105").Value = 1121
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 262
$ws.Range("E105").Value = 849
$ws.Range("H105").Value = 10

$ws.Range("A106").Value = 'Libano'
$ws.Range("B106").Value = 1119
$ws.Range("C106").Value = 5
$ws.Range("D106").Value = 688
$ws.Range("E106").Value = 405
$ws.Range("H106").Value = 26

# --- Rows 114-116: Niger / Republica de Chipre / Costa Rica rotate ---
$ws.Range("A114").Value = 'Costa Rica'
$ws.Range("B114").Value = 951
$ws.Range("C114").Value = 21
$ws.Range("D114").Value = 628
$ws.Range("E114").Value = 313
$ws.Range("H114").Value = 10

$ws.Range("A115").Value = 'Niger'
$ws.Range("B115").Value = 945
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 783
$ws.Range("E115").Value = 101
$ws.Range("H115").Value = 61

$ws.Range("A116").Value = 'Republica de Chipre'
$ws.Range("B116").Value = 937
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 594
$ws.Range("E116").Value = 326
$ws.Range("H116").Value = 17

# --- Rows 151-152: Suazilandia <-> Mauritania swap ---
$ws.Range("A151").Value = 'Mauritania'
$ws.Range("B151").Value = 262
$ws.Range("C151").Value = 25
$ws.Range("D151").Value = 15
$ws.Range("E151").Value = 238
$ws.Range("G151").Value = 3
$ws.Range("H151").Value = 9

$ws.Range("A152").Value = 'Suazilandia'
$ws.Range("B152").Value = 256
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 158
$ws.Range("E152").Value = 96
$ws.Range("H152").Value = 2

# --- Rows 197-198: Fiyi <-> Curazao swap ---
$ws.Range("A197").Value = 'Curazao'
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 14
$ws.Range("H197").Value = 1

$ws.Range("A198").Value = 'Fiyi'
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

# --- Rows 199-201: Santa Lucia / Nueva Caledonia / Belice rotate ---
$ws.Range("A199").Value = 'Nueva Caledonia'

$ws.Range("A200").Value = 'Belice'
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = 'Santa Lucia'
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# --- Rows 207-208: Islas Turcas y Caicos <-> Groenlandia swap ---
$ws.Range("A207").Value = 'Groenlandia'
$ws.Range("D207").Value = 11
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = 'Islas Turcas y Caicos'
$ws.Range("D208").Value = 10
$ws.Range("H208").Value = 1

# --- Rows 210-211: Montserrat <-> Seychelles swap ---
$ws.Range("A210").Value = 'Seychelles'
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = 'Montserrat'
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
